# Append 45 new data rows (rows 102-146) to the "master-reg_center_machine_devic"
# sheet, following the same repeating pattern already present in the sheet
# (A/B cycle through a fixed block of 9 values, C increments by 1, D/E/F/G
# are constant lookup/boilerplate values), then update the visible
# selection/view and the page orientation to match the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseA = 10002
$baseB = 10021
$baseC = 3000121

for ($i = 0; $i -lt 45; $i++) {
    $r = 102 + $i
    $ws.Cells.Item($r, 1).Value = $baseA + ($i % 9)
    $ws.Cells.Item($r, 2).Value = $baseB + ($i % 9)
    $ws.Cells.Item($r, 3).Value = $baseC + $i
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Match the saved view state: the new block (A102:G146) is selected, with
# A102 as the active cell.
$ws.Range("A102:G146").Select()

# The workbook was re-saved with the print orientation set to Portrait.
$ws.PageSetup.Orientation = 1
